$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.659.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.438.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.38%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.65%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.435.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("E10").Value = "  -5.48%  "
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("E13").Value = "  -3.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000174"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.893.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.442.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.441.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.16%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "637.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.44%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0961"
$ws.Range("E28").Value = "  -9.28%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.565.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -6.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.139"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("E37").Value = "  -6.01%  "
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0524"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.598"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.50%  "
$ws.Range("D51").Value = "0.0₆0235"
$ws.Range("E51").Value = "  +6.66%  "
